$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Lockie Ferguson"
$ws.Name = "Lockie Ferguson"

# Header row (insert new "matchNo" column at A, shifting everything right)
$ws.Range("A1").Value = "'matchNo"
$ws.Range("B1").Value = "'teamName"
$ws.Range("C1").Value = "'batterName"
$ws.Range("D1").Value = "'states"
$ws.Range("E1").Value = "'runs"
$ws.Range("F1").Value = "'balls"
$ws.Range("G1").Value = "'fours"
$ws.Range("H1").Value = "'sixes"
$ws.Range("I1").Value = "'sr"
$ws.Range("J1").Value = "'opponentTeamName"
$ws.Range("K1").Value = "'venue"
$ws.Range("L1").Value = "'date"
$ws.Range("M1").Value = "'result"

# Row 2: new "Qualifier" match row
$ws.Range("A2").Value = "'Qualifier"
$ws.Range("B2").Value = "'Kolkata Knight Riders"
$ws.Range("C2").Value = "'Lockie Ferguson"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'0"
$ws.Range("F2").Value = "'0"
$ws.Range("G2").Value = "'0"
$ws.Range("H2").Value = "'0"
$ws.Range("I2").Value = "'-"
$ws.Range("J2").Value = "'Delhi Capitals"
$ws.Range("K2").Value = "'Sharjah"
$ws.Range("L2").Value = "'October 13"
$ws.Range("M2").Value = "'KKR won by 3 wickets (with 1 ball remaining)"

# Row 3: new "41st" match row
$ws.Range("A3").Value = "'41st"
$ws.Range("B3").Value = "'Kolkata Knight Riders"
$ws.Range("C3").Value = "'Lockie Ferguson"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = "'0"
$ws.Range("F3").Value = "'0"
$ws.Range("G3").Value = "'0"
$ws.Range("H3").Value = "'0"
$ws.Range("I3").Value = "'-"
$ws.Range("J3").Value = "'Delhi Capitals"
$ws.Range("K3").Value = "'Sharjah"
$ws.Range("L3").Value = "'September 28"
$ws.Range("M3").Value = "'KKR won by 3 wickets (with 10 balls remaining)"

# Row 4: the "Final" match row (originally row 2, now shifted down + right by one col)
$ws.Range("A4").Value = "'Final"
$ws.Range("B4").Value = "'Kolkata Knight Riders"
$ws.Range("C4").Value = "'Lockie Ferguson"
$ws.Range("D4").Value = "'"
$ws.Range("E4").Value = "'18"
$ws.Range("F4").Value = "'11"
$ws.Range("G4").Value = "'1"
$ws.Range("H4").Value = "'1"
$ws.Range("I4").Value = "'163.63"
$ws.Range("J4").Value = "'Chennai Super Kings"
$ws.Range("K4").Value = "'Dubai (DSC)"
$ws.Range("L4").Value = "'October 15"
$ws.Range("M4").Value = "'Super Kings won by 27 runs"
